$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the refreshed cryptocurrency data feed
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.198.47'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.678.19'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.55'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.64%  '
$ws.Range('E9').Value = '  +3.10%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.915.49'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.20'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.655.64'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.560'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.70'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.153.38'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '235.61'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.85'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.55'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.55'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.15%  '
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.77'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.43'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.537.40'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('E35').Value = '  -3.20%  '
$ws.Range('E36').Value = '  +3.71%  '
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('E40').Value = '  +2.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '70.02'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.19%  '
$ws.Range('E42').Value = '  +4.65%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.823.35'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('E47').Value = '  +7.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.97'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.36%  '
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.24'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.08%  '
$ws.Range('E51').Value = '  +1.59%  '
